# Pedidos.xlsx update:
# Move the 13 newest order rows (Remessas 80266975-80266930) from the
# bottom of the data block up to the top of the list (right after the
# header row), shifting the existing data down. The vacated rows at the
# old location are left blank. Also park the selection at D5 (top of the
# newly inserted block) instead of the old D217.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 13 blank rows above row 2 (the first data row), pushing all
#    existing data (and the 13 "new" rows that currently sit at the foot
#    of the table, rows 218-230) down by 13 rows.
$ws.Range("A2:A14").EntireRow.Insert()

# 2) The 13 newest orders, which used to be rows 218-230, are now sitting
#    (still fully populated) at rows 231-243. Fill the freshly inserted
#    rows 2-14 with that same data.
$newOrders = @(
  @(80266975, "10661-BLB-I", 10),
  @(80266975, "10661-BLB-I", 10),
  @(80266976, "10636-ARI-I", 2),
  @(80266976, "10636-ARI-I", 2),
  @(80266977, "33664-ATE-I", 1),
  @(80266977, "10251-ARI-I", 1),
  @(80266977, "21475-NZX-I", 1),
  @(80266977, "10080-ARI-I", 4),
  @(80266977, "23364-ATE-I", 4),
  @(80266978, "10661-BLB-I", 11),
  @(80266983, "10662-BLB-I", 58),
  @(80266930, "10663-BLB-I", 8),
  @(80266930, "10072-BLB-I", 5)
)

for ($i = 0; $i -lt $newOrders.Count; $i++) {
  $targetRow = 2 + $i
  $ws.Cells.Item($targetRow, 1).Value = $newOrders[$i][0]
  $ws.Cells.Item($targetRow, 2).Value = $newOrders[$i][1]
  $ws.Cells.Item($targetRow, 3).Value = $newOrders[$i][2]
}

# 3) Clear out the now-duplicated data left behind at rows 231-243 (their
#    old location before the insert), leaving blank cells in place.
$ws.Range("A231:C243").ClearContents()

# 4) Move the visible selection to D5 (previously D217).
[void]$ws.Range("D5").Select()
